# Regenerate save_data "K" column (formerly "Strike#") for each row, per
# the updated std/mean calculation in the data-regeneration pipeline.
#
# Column map on Sheet1: B=date C=TB D=PC E=dS0 F=dSF G=K H=IP I=I0 J=IF
# This commit rewrites the computed K values (column G) for rows 2-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 3
    30 = 0
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 1
    48 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
